# "tablas y evaluaciones todas"
# Remove the now-unused trailing column (V) from the repeat_p24_4 sheet
# (its only used cell was the header "p24_4" in V1) and reset the
# data column U to 0 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column V entirely (header "p24_4" + no data cells below it).
$ws.Columns.Item(22).Delete()

# Column U (now the last column) holds 0 for every data row 2-69,
# including row 59 which previously held a stray 1.6 value.
$ws.Range("U2:U69").Value = 0
